$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 2): S2:W2 ---
$ws.Range("S2").Value = "Related Item"
$ws.Range("T2").Value = "Related Item Label"
$ws.Range("U2").Value = "Terms of Use"
$ws.Range("V2").Value = "Language"
$ws.Range("W2").Value = "Original Physical Description"

# --- Row 3 updates: R3 changes from "Catalog Key" to "local", plus new S3:W3 ---
$ws.Range("R3").Value = "local"
$ws.Range("S3").Value = "http://www.example.com/text.pdf"
$ws.Range("T3").Value = "Example Item PDF"
$ws.Range("U3").Value = "Terms of Use Language"
$ws.Range("V3").Value = "English"
$ws.Range("W3").Value = "16mm Reel"

# --- Column widths for the newly-introduced columns (Q=17 .. W=23) ---
# The runtime's ColumnWidth setter rounds to whole pixels at MDW=6, with a
# fixed +5/6 offset baked into the OOXML width it writes back out, so the
# inputs below are pre-compensated to land as close as possible on the
# target OOXML widths (8.83333333333333, 12.537037037037, 27.6666666666667,
# 16.2666666666667, 19.9962962962963, 8.83333333333333, 23.937037037037).
$ws.Columns.Item(17).ColumnWidth = 8.0
$ws.Columns.Item(18).ColumnWidth = 11.666666666666666
$ws.Columns.Item(19).ColumnWidth = 26.833333333333332
$ws.Columns.Item(20).ColumnWidth = 15.5
$ws.Columns.Item(21).ColumnWidth = 19.166666666666668
$ws.Columns.Item(22).ColumnWidth = 8.0
$ws.Columns.Item(23).ColumnWidth = 23.166666666666668

# --- Selection / active cell moves from R3 to U2 (matches the target sheetView) ---
$ws.Range("U2").Select()
